# Weekly update to the "Hortaliza, Terminal Hortofrutícola Agro Chillán - Papa" sheet.
# A new price observation is inserted as row 325 (pushing the existing rows
# 325-369 down to 326-370), adding a fresh "1a (guarda)" Patagonia potato
# quote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 325, shifting rows 325:369 down to 326:370.
$ws.Rows.Item(325).Insert()

# Populate the newly inserted row 325 with the new weekly observation.
$ws.Cells.Item(325, 1).Value  = 7
$ws.Cells.Item(325, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(325, 3).Value  = "Ñuble"
$ws.Cells.Item(325, 4).Value  = 44776
$ws.Cells.Item(325, 5).Value  = 16
$ws.Cells.Item(325, 6).Value  = 100114001
$ws.Cells.Item(325, 7).Value  = "Papa"
$ws.Cells.Item(325, 8).Value  = "Patagonia"
$ws.Cells.Item(325, 9).Value  = "1a (guarda)"
$ws.Cells.Item(325, 10).Value = 160
$ws.Cells.Item(325, 11).Value = 7000
$ws.Cells.Item(325, 12).Value = 7500
$ws.Cells.Item(325, 13).Value = 7250
$ws.Cells.Item(325, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(325, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(325, 16).Value = 290
$ws.Cells.Item(325, 17).Value = 25
$ws.Cells.Item(325, 18).Value = "Hortaliza"
